$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '30.397.38'
$ws.Range('E2').Value = '  -1.75%  '
$ws.Range('D3').Value = '1.915.62'
$ws.Range('E3').Value = '  +1.28%  '
Set-TextValue $ws 'D4' '0.9995'
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue $ws 'D5' '241.24'
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('E6').Value = '  -0.07%  '
Set-TextValue $ws 'D7' '0.4694'
$ws.Range('E7').Value = '  -2.11%  '
Set-TextValue $ws 'D8' '0.2841'
$ws.Range('E8').Value = '  -0.58%  '
Set-TextValue $ws 'D9' '0.06886'
$ws.Range('E9').Value = '  +5.15%  '
Set-TextValue $ws 'D10' '106.67'
$ws.Range('E10').Value = '  +10.95%  '
Set-TextValue $ws 'D11' '18.05'
$ws.Range('E11').Value = '  -3.92%  '
$ws.Range('D12').Value = '1.906.92'
$ws.Range('E12').Value = '  +0.74%  '
Set-TextValue $ws 'D13' '0.07635'
$ws.Range('E13').Value = '  +1.32%  '
Set-TextValue $ws 'D14' '5.189'
$ws.Range('E14').Value = '  +0.96%  '
Set-TextValue $ws 'D15' '0.6548'
$ws.Range('E15').Value = '  +0.16%  '
Set-TextValue $ws 'D16' '287.42'
$ws.Range('E16').Value = '  -3.68%  '
$ws.Range('D17').Value = '30.399.90'
$ws.Range('E17').Value = '  -1.62%  '
$ws.Range('E18').Value = '  -1.54%  '
Set-TextValue $ws 'D19' '1.000'
$ws.Range('E19').Value = '  +0.39%  '
Set-TextValue $ws 'D20' '0.000007590'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('D21').Value = '2.145.38'
$ws.Range('E21').Value = '  +1.17%  '
Set-TextValue $ws 'D22' '0.9996'
$ws.Range('E22').Value = '  +0.07%  '
Set-TextValue $ws 'D23' '5.214'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('E24').Value = '  +0.56%  '
Set-TextValue $ws 'D25' '168.10'
$ws.Range('E25').Value = '  -0.10%  '
Set-TextValue $ws 'D26' '9.240'
$ws.Range('E26').Value = '  -0.72%  '
Set-TextValue $ws 'D27' '21.20'
$ws.Range('E27').Value = '  +6.84%  '
Set-TextValue $ws 'D28' '2.031'
$ws.Range('E28').Value = '  +2.93%  '
Set-TextValue $ws 'D29' '0.1073'
$ws.Range('E29').Value = '  +1.40%  '
Set-TextValue $ws 'D30' '1.367'
$ws.Range('E30').Value = '  -1.49%  '
Set-TextValue $ws 'D31' '4.133'
$ws.Range('E31').Value = '  -0.87%  '
Set-TextValue $ws 'D32' '3.945'
$ws.Range('E32').Value = '  -1.04%  '
Set-TextValue $ws 'D33' '0.05043'
$ws.Range('E33').Value = '  +0.45%  '
Set-TextValue $ws 'D34' '0.7357'
$ws.Range('E34').Value = '  +1.22%  '
Set-TextValue $ws 'D35' '1.142'
$ws.Range('E35').Value = '  -3.32%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws 'D36' '2.721'
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D37' '0.02002'
$ws.Range('E37').Value = '  +3.15%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D38' '2.668'
$ws.Range('E38').Value = '  -1.79%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D39' '2.045'
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D40' '108.71'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws 'D41' '0.8731'
$ws.Range('E41').Value = '  -3.03%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D42' '5.831'
$ws.Range('E42').Value = '  +3.95%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws 'D43' '0.9998'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws 'D44' '52.62'
$ws.Range('E44').Value = '  +24.99%  '
$ws.Range('E45').Value = '  -0.83%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D46' '67.08'
$ws.Range('E46').Value = '  +2.32%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws 'D47' '7.111'
$ws.Range('E47').Value = '  -3.75%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D48' '9.151'
$ws.Range('E48').Value = '  +2.51%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D49' '0.1202'
$ws.Range('E49').Value = '  -2.35%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws 'D50' '34.56'
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('B51').Value = 'eCash'
$ws.Range('C51').Value = 'https://coinranking.com/coin/aQx_vW8s1+ecash-xec'
Set-TextValue $ws 'D51' '0.00004446'
$ws.Range('E51').Value = '  +57.45%  '
